# Update the "6.Adjustments" sheet: replace the "No impact" text entries
# with numeric 0 values, and make every data cell in B2:M10 use the plain
# two-decimal numeric format instead of the old percentage style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6.Adjustments")
$ws.Activate()

# Row -> values for columns B..M (12 values each, col B = 2 .. col M = 13)
$rowValues = @{
    2  = @(0, 0.25, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5)
    3  = @(0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
    4  = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    5  = @(0, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75)
    6  = @(0, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75)
    7  = @(0, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75, 0.75)
    8  = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    9  = @(0, 0.25, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5)
    10 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}

# All data cells now share the plain 2-decimal numeric format (was a mix of
# percentage-format "No impact" text cells and 2-decimal numeric cells).
$ws.Range("B2:M10").NumberFormat = "0.00"

# Restore selection to B4:M4 with active cell B4 (matches the saved view state)
$ws.Range("B4:M4").Select() | Out-Null
